# Se añadió la tabla de tiempos
# Adds a new worksheet "Hoja2" after "Hoja1" containing a summary table of
# times, pulling values from Hoja1 via formulas.

$wb = $excel.ActiveWorkbook
$hoja1 = $wb.Worksheets.Item("Hoja1")

# Make sure Hoja1 is the active sheet so the new sheet is inserted right
# after it.
$hoja1.Activate()

$hoja2 = $wb.Worksheets.Add($null, $hoja1)
$hoja2.Name = "Hoja2"

# Header row (bold, no special alignment)
$hoja2.Range("B2").Value = "Paradas"
$hoja2.Range("C2").Value = "SECUENCIAL"
$hoja2.Range("D2").Value = "AUTO"
$hoja2.Range("E2").Value = "DYNAMIC"
$hoja2.Range("F2").Value = "GUIDED"
$hoja2.Range("G2").Value = "STATIC"
$hoja2.Range("B2:G2").Font.Bold = $true

# Data rows: column B holds the "number of stops" literal, columns C..G
# pull the corresponding timing figures from Hoja1 via formulas.
$rows = @(
    @{ Row2 = 3;  Stops = 2;  Row1 = 5 },
    @{ Row2 = 4;  Stops = 4;  Row1 = 6 },
    @{ Row2 = 5;  Stops = 5;  Row1 = 7 },
    @{ Row2 = 6;  Stops = 6;  Row1 = 8 },
    @{ Row2 = 7;  Stops = 7;  Row1 = 9 },
    @{ Row2 = 8;  Stops = 8;  Row1 = 10 },
    @{ Row2 = 9;  Stops = 9;  Row1 = 11 },
    @{ Row2 = 10; Stops = 10; Row1 = 12 }
)

foreach ($r in $rows) {
    $r2 = $r.Row2
    $r1 = $r.Row1
    $hoja2.Range("B$r2").Value = $r.Stops
    $hoja2.Range("C$r2").Formula = "=Hoja1!I$r1"
    $hoja2.Range("D$r2").Formula = "=Hoja1!H$r1"
    $hoja2.Range("E$r2").Formula = "=Hoja1!C$r1"
    $hoja2.Range("F$r2").Formula = "=Hoja1!M$r1"
    $hoja2.Range("G$r2").Formula = "=Hoja1!R$r1"
}

# Column C is a bit wider to fit the "SECUENCIAL" header.
$hoja2.Range("C1").ColumnWidth = 13.29

# Keep the same portrait page setup used by Hoja1.
$hoja2.PageSetup.Orientation = 1

# Match the saved view state: Hoja1 keeps its Q31 selection but scrolls
# back so column F is visible at the left edge, while Hoja2 ends up as the
# active/selected sheet with G5 selected.
# NOTE: Range.Select() activates its parent sheet, so Hoja1's window must
# be configured first and Hoja2 activated last so it remains the active tab.
$hoja1.Activate()
$win1 = $excel.ActiveWindow
$win1.ScrollColumn = 6
$win1.ScrollRow = 1

$hoja2.Activate()
[void]$hoja2.Range("G5").Select()
